$p = $ppt.ActivePresentation
$p.Slides.Item(5).Delete()
